# Sprint3Backlog.xlsx - update Priority for "Testing for boards page" (row 11)
# from 1 (Nate) to 2 (Nick), matching the added Xunit login-controller testing work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Range("B11")
$cell.Value = 2
$cell.Select()
